$wb = $excel.ActiveWorkbook

$wsParameters = $wb.Worksheets.Item("Parameters")
$wsAlgorithm  = $wb.Worksheets.Item("Algorithm")
$wsOperators  = $wb.Worksheets.Item("Operators")

# ---------------------------------------------------------------------------
# The order in which brand-new strings are first written determines their
# position in the shared-strings table, so the very first pass below touches
# only the cells that introduce a NEW string, in the exact sequence required
# (mj_zdt1_decimal, tools, selNSGA2, mj_evaluators, mj_algorithms, int,
# value, type, mutate, mj_operators, selTournamentDCD, Probability crossover,
# mate, Probability flip allele, Jump size, Probability mutation).
# ---------------------------------------------------------------------------

$wsAlgorithm.Range("C2").Value = "mj_zdt1_decimal"
$wsAlgorithm.Range("B3").Value = "tools"
$wsAlgorithm.Range("C3").Value = "selNSGA2"
$wsAlgorithm.Range("B2").Value = "mj_evaluators"
$wsAlgorithm.Range("B1").Value = "mj_algorithms"

$wsParameters.Range("C2").Value = "int"
$wsParameters.Range("B1").Value = "value"
$wsParameters.Range("C1").Value = "type"

$wsOperators.Range("A3").Value = "mutate"
$wsOperators.Range("B2").Value = "mj_operators"
$wsOperators.Range("C1").Value = "selTournamentDCD"

$wsParameters.Range("A4").Value = "Probability crossover"

$wsOperators.Range("A2").Value = "mate"

$wsParameters.Range("A5").Value = "Probability flip allele"
$wsParameters.Range("A6").Value = "Jump size"
$wsParameters.Range("A7").Value = "Probability mutation"

# ---------------------------------------------------------------------------
# Parameters sheet - fill in the rest of the table (reuses strings above).
# ---------------------------------------------------------------------------
$wsParameters.Range("A1").Value = "name"

$wsParameters.Range("A2").Value = "Population size"
$wsParameters.Range("B2").Value = 40

$wsParameters.Range("A3").Value = "Generations"
$wsParameters.Range("B3").Value = 10
$wsParameters.Range("C3").Value = "int"

$wsParameters.Range("B4").Value = 0.5
$wsParameters.Range("C4").Value = "float"

$wsParameters.Range("B5").Value = 0.5
$wsParameters.Range("C5").Value = "float"

$wsParameters.Range("B6").Value = 10
$wsParameters.Range("C6").Value = "int"

$wsParameters.Range("B7").Value = 0.5
$wsParameters.Range("C7").Value = "float"

[void]$wsParameters.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# Algorithm sheet - fill in the rest of the table.
# ---------------------------------------------------------------------------
$wsAlgorithm.Range("A1").Value = "name"
$wsAlgorithm.Range("C1").Value = "nsga2"

$wsAlgorithm.Range("A2").Value = "evaluate"

$wsAlgorithm.Range("A3").Value = "select"

[void]$wsAlgorithm.Columns.Item(2).AutoFit()
$wsAlgorithm.PageSetup.PaperSize = 9
$wsAlgorithm.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Operators sheet - fill in the rest of the table.
# ---------------------------------------------------------------------------
$wsOperators.Range("A1").Value = "select"
$wsOperators.Range("A1").Font.Bold = $false
$wsOperators.Range("B1").Value = "tools"

$wsOperators.Range("C2").Value = "mj_list_flip"

$wsOperators.Range("B3").Value = "mj_operators"
$wsOperators.Range("C3").Value = "mj_random_jump"

[void]$wsOperators.Columns.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# Selections / active sheet.
# ---------------------------------------------------------------------------
[void]$wsAlgorithm.Range("C9").Select()
[void]$wsOperators.Range("A3").Select()

[void]$wsParameters.Activate()
[void]$wsParameters.Range("C7").Select()
